$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference style (default body-row style) used to restore formatting
# after forcing a Text number format so numeric-looking strings are not
# auto-converted to the Number type (matches original inlineStr cells).
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "27.201.60"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "1.820.20"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  -1.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.92"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4258"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.99"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07207"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8600"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.96"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "1.838.56"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.649"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07098"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.300"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  -3.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.72"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008860"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.01"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "27.238.36"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.114"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "2.070.53"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.006"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.91"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.108"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +5.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.222"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.11"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08870"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7589"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.190"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.459"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.826"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -6.84%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.114"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05242"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.900"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.040"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1674"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5009"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.619"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.58"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.23"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06405"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.657"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -3.34%  "
